$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Colour the "Importing product data from MongoDB into SQL Server"
#    bullet green (RGB 00B050). WdColor values are stored as BGR, so
#    00B050 (RGB) -> 50B000 (BGR) -> 5287936 decimal.
# -----------------------------------------------------------------
$mongoRange = $d.Content
$mongoFound = $mongoRange.Find.Execute(
    "Importing product data from MongoDB into SQL Server",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($mongoFound) {
    # $mongoRange.Paragraphs.Item(1).Range is clipped to the searched
    # text and excludes the paragraph mark, so the mark's rPr would be
    # skipped. Locate the real Document paragraph instead, which
    # covers the mark too (diff colours it as well).
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $mongoRange.Start -and $mongoRange.Start -lt $p.Range.End) {
            $p.Range.Font.Color = 5287936
            break
        }
    }
}

# -----------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so that it now sits right after the
#    "(ADO.NET)" run (it currently sits inside "Importing data from
#    XML"). Adding a bookmark with the same name automatically drops
#    the previous one, which is exactly what the diff shows.
#
#    NOTE: this runtime mis-places a *collapsed* Range/bookmark whose
#    Start lands exactly on a paragraph-end (the character right
#    before the paragraph mark) - it resets to the start of the
#    document instead. To dodge that, we temporarily insert a marker
#    character right after "(ADO.NET)", anchor the bookmark next to
#    it (now a safe, non-paragraph-end position), and then delete the
#    marker again; the bookmark stays correctly anchored.
# -----------------------------------------------------------------
$adoRange = $d.Content
$adoFound = $adoRange.Find.Execute("(ADO.NET)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($adoFound) {
    $afterAdo = $adoRange.End

    $marker = $d.Range($afterAdo, $afterAdo)
    $marker.InsertAfter("@")

    $bmAnchor = $d.Range($afterAdo, $afterAdo)
    $d.Bookmarks.Add("_GoBack", $bmAnchor)

    $markerRange = $d.Range($afterAdo, $afterAdo + 1)
    $markerRange.Delete()
}

# -----------------------------------------------------------------
# 3) Collapse the two runs that used to straddle the old bookmark
#    ("Importing data fro" / "m XML") into a single run reading
#    "Importing data from XML".
# -----------------------------------------------------------------
$xmlRange = $d.Content
$xmlRange.Find.Execute(
    "Importing data from XML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Importing data from XML", 2) | Out-Null
